# "Fruta / hortaliza, semanal"
#
# The weekly refresh adds six new daily price records for Damasco (Vega
# Central Mapocho de Santiago): two at the very top of the existing block
# (rows 28-29, dated 2021-12-24) and four more further down (rows 45-48,
# dated 2021-12-23), pushing every subsequent record down accordingly.
# Final used range grows from A1:T60 to A1:T66.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Make room for the new rows -----------------------------------
# Two blank rows right before the old row 28 (old 28-60 -> new 30-62).
$ws.Rows.Item(28).Insert()
$ws.Rows.Item(28).Insert()

# Four more blank rows before what is now row 45 (old 43-60, already
# shifted to 45-62, move again to 49-66).
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()

# --- 2) Values shared by every row of this product block -------------
$common = @{
    1  = 9
    2  = "Vega Central Mapocho de Santiago"
    3  = "Metropolitana"
    5  = 13
    6  = "Fruta"
    7  = 100103
    8  = "Frutos de hueso (carozo)"
    9  = 100103003
    10 = "Damasco"
}

# --- 3) Per-row specifics for the six brand-new records ---------------
# columns: D, K, L, M, N, O, P, Q, R, S, T
$newRows = @{
    28 = @{ D = 44554; K = "Patterson"; L = "Primera"; M = 410; N = 14000; O = 14000; P = 14000; Q = "`$/caja 18 kilos";    R = "Provincia de San Felipe de Aconcagua"; S = 778;  T = 18 }
    29 = @{ D = 44554; K = "Patterson"; L = "Segunda"; M = 380; N = 12000; O = 12000; P = 12000; Q = "`$/caja 18 kilos";    R = "Provincia de San Felipe de Aconcagua"; S = 667;  T = 18 }
    45 = @{ D = 44553; K = "Dina";      L = "Especial"; M = 280; N = 7500;  O = 7500;  P = 7500;  Q = "`$/bandeja 6 kilos"; R = "Provincia de San Felipe de Aconcagua"; S = 1250; T = 6  }
    46 = @{ D = 44553; K = "Dina";      L = "Primera";  M = 350; N = 6500;  O = 6500;  P = 6500;  Q = "`$/bandeja 6 kilos"; R = "Provincia de San Felipe de Aconcagua"; S = 1083; T = 6  }
    47 = @{ D = 44553; K = "Patterson"; L = "Primera";  M = 380; N = 14000; O = 14000; P = 14000; Q = "`$/caja 18 kilos";   R = "Provincia de San Felipe de Aconcagua"; S = 778;  T = 18 }
    48 = @{ D = 44553; K = "Patterson"; L = "Segunda";  M = 410; N = 12000; O = 12000; P = 12000; Q = "`$/caja 18 kilos";   R = "Provincia de San Felipe de Aconcagua"; S = 667;  T = 18 }
}

$colLetterToIndex = @{ D = 4; K = 11; L = 12; M = 13; N = 14; O = 15; P = 16; Q = 17; R = 18; S = 19; T = 20 }

foreach ($rowNum in 28, 29, 45, 46, 47, 48) {
    foreach ($col in $common.Keys) {
        $ws.Cells.Item($rowNum, $col).Value = $common[$col]
    }
    $rowData = $newRows[$rowNum]
    foreach ($letter in $rowData.Keys) {
        $colIndex = $colLetterToIndex[$letter]
        $ws.Cells.Item($rowNum, $colIndex).Value = $rowData[$letter]
    }
}

Write-Output "Inserted 6 new Damasco price rows; used range now $($ws.Range('A1').Worksheet.UsedRange.Address())"
